$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on column D so numeric-looking strings keep their exact text representation
# (mirrors the original workbook where D/E are stored as inline strings, not numbers)
$ws.Range("D2:D51").NumberFormat = "@"

# Rows 34 and 35 swap content (ImmutableX <-> WEMIXTOKEN) and get new price/volume figures
$ws.Cells.Item(34, 2).Value = "WEMIXTOKEN"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(34, 4).Value = "2.008"
$ws.Cells.Item(34, 5).Value = "  -5.89%  "

$ws.Cells.Item(35, 2).Value = "ImmutableX"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(35, 4).Value = "1.001"
$ws.Cells.Item(35, 5).Value = "  +1.26%  "

# Price/Volume updates for the remaining rows
$ws.Cells.Item(2, 4).Value = '22.355.96'
$ws.Cells.Item(2, 5).Value = '  -4.64%  '
$ws.Cells.Item(3, 4).Value = '1.562.87'
$ws.Cells.Item(3, 5).Value = '  -5.09%  '
$ws.Cells.Item(4, 5).Value = '  +0.39%  '
$ws.Cells.Item(5, 4).Value = '1.002'
$ws.Cells.Item(5, 5).Value = '  +0.30%  '
$ws.Cells.Item(6, 4).Value = '289.62'
$ws.Cells.Item(6, 5).Value = '  -3.29%  '
$ws.Cells.Item(7, 4).Value = '0.3698'
$ws.Cells.Item(7, 5).Value = '  -2.57%  '
$ws.Cells.Item(8, 5).Value = '  -2.30%  '
$ws.Cells.Item(9, 4).Value = '0.3375'
$ws.Cells.Item(9, 5).Value = '  -3.51%  '
$ws.Cells.Item(10, 4).Value = '1.162'
$ws.Cells.Item(10, 5).Value = '  -4.76%  '
$ws.Cells.Item(11, 4).Value = '0.07637'
$ws.Cells.Item(11, 5).Value = '  -5.44%  '
$ws.Cells.Item(12, 5).Value = '  +0.40%  '
$ws.Cells.Item(13, 4).Value = '21.32'
$ws.Cells.Item(13, 5).Value = '  -3.47%  '
$ws.Cells.Item(14, 4).Value = '6.040'
$ws.Cells.Item(14, 5).Value = '  -4.37%  '
$ws.Cells.Item(15, 4).Value = '6.910'
$ws.Cells.Item(15, 5).Value = '  -5.10%  '
$ws.Cells.Item(16, 4).Value = '1.558.41'
$ws.Cells.Item(16, 5).Value = '  -4.57%  '
$ws.Cells.Item(17, 4).Value = '0.00001125'
$ws.Cells.Item(17, 5).Value = '  -7.78%  '
$ws.Cells.Item(18, 4).Value = '90.11'
$ws.Cells.Item(18, 5).Value = '  -5.15%  '
$ws.Cells.Item(19, 4).Value = '0.06730'
$ws.Cells.Item(19, 5).Value = '  -3.39%  '
$ws.Cells.Item(20, 5).Value = '  +0.31%  '
$ws.Cells.Item(21, 4).Value = '6.248'
$ws.Cells.Item(21, 5).Value = '  -5.78%  '
$ws.Cells.Item(22, 4).Value = '16.52'
$ws.Cells.Item(22, 5).Value = '  -4.90%  '
$ws.Cells.Item(23, 4).Value = '0.5273'
$ws.Cells.Item(23, 5).Value = '  -7.30%  '
$ws.Cells.Item(24, 4).Value = '12.01'
$ws.Cells.Item(24, 5).Value = '  -3.40%  '
$ws.Cells.Item(25, 4).Value = '22.346.01'
$ws.Cells.Item(25, 5).Value = '  -4.68%  '
$ws.Cells.Item(26, 4).Value = '2.371'
$ws.Cells.Item(26, 5).Value = '  -2.73%  '
$ws.Cells.Item(27, 4).Value = '2.789'
$ws.Cells.Item(27, 5).Value = '  -6.56%  '
$ws.Cells.Item(28, 5).Value = '  -4.30%  '
$ws.Cells.Item(29, 4).Value = '145.91'
$ws.Cells.Item(29, 5).Value = '  -2.46%  '
$ws.Cells.Item(30, 4).Value = '4.974'
$ws.Cells.Item(30, 5).Value = '  -3.95%  '
$ws.Cells.Item(31, 4).Value = '125.49'
$ws.Cells.Item(31, 5).Value = '  -4.73%  '
$ws.Cells.Item(32, 4).Value = '1.732.98'
$ws.Cells.Item(32, 5).Value = '  -4.78%  '
$ws.Cells.Item(33, 4).Value = '6.200'
$ws.Cells.Item(33, 5).Value = '  -9.25%  '
$ws.Cells.Item(36, 4).Value = '10.04'
$ws.Cells.Item(36, 5).Value = '  -10.95%  '
$ws.Cells.Item(37, 4).Value = '0.08433'
$ws.Cells.Item(37, 5).Value = '  -3.92%  '
$ws.Cells.Item(38, 5).Value = '  -5.85%  '
$ws.Cells.Item(39, 4).Value = '0.2316'
$ws.Cells.Item(39, 5).Value = '  -4.58%  '
$ws.Cells.Item(40, 4).Value = '5.515'
$ws.Cells.Item(40, 5).Value = '  -6.37%  '
$ws.Cells.Item(41, 4).Value = '0.06439'
$ws.Cells.Item(41, 5).Value = '  -5.73%  '
$ws.Cells.Item(42, 4).Value = '1.288'
$ws.Cells.Item(42, 5).Value = '  -0.38%  '
$ws.Cells.Item(43, 4).Value = '11.65'
$ws.Cells.Item(43, 5).Value = '  -8.90%  '
$ws.Cells.Item(44, 4).Value = '0.6319'
$ws.Cells.Item(44, 5).Value = '  -7.62%  '
$ws.Cells.Item(45, 4).Value = '14.18'
$ws.Cells.Item(45, 5).Value = '  -9.13%  '
$ws.Cells.Item(46, 5).Value = '  +0.31%  '
$ws.Cells.Item(47, 4).Value = '0.5965'
$ws.Cells.Item(47, 5).Value = '  -6.03%  '
$ws.Cells.Item(48, 4).Value = '3.754'
$ws.Cells.Item(48, 5).Value = '  -4.04%  '
$ws.Cells.Item(49, 4).Value = '2.095'
$ws.Cells.Item(49, 5).Value = '  -6.56%  '
$ws.Cells.Item(50, 4).Value = '1.261'
$ws.Cells.Item(50, 5).Value = '  +3.43%  '
$ws.Cells.Item(51, 4).Value = '123.85'
$ws.Cells.Item(51, 5).Value = '  -2.58%  '
